# "Generate Report for Handoff"
# Refresh the handoff-tracking report: the source markdown file was re-handed-off
# under a new GUID/commit hash, so update the file name references, the
# generated .xlf target file names, and the handoff timestamps across all
# three report sheets (Overview, zh-cn, de-de), keeping hyperlink display text
# in sync with the cell values.

$wb = $excel.ActiveWorkbook

$newGuid = "ca5a85e7-5afe-49ec-9f16-90241ba6d507"
$newHash = "c76724436579235c4a6c9419922143b25703047b"

$newMdName    = "$newGuid.md"
$newZhXlfName = "$newGuid.$newHash.zh-cn.xlf"
$newDeXlfName = "$newGuid.$newHash.de-de.xlf"

function Set-HyperlinkDisplay($ws, $row, $col, $text) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Row -eq $row -and $hl.Range.Column -eq $col) {
            $hl.TextToDisplay = $text
        }
    }
}

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newMdName
Set-HyperlinkDisplay $wsOverview 2 1 $newMdName

$wsOverview.Range("D2").Value = "2016-03-13 09:03:33"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = $newMdName
Set-HyperlinkDisplay $wsZhCn 2 1 $newMdName

$wsZhCn.Range("D2").Value = $newZhXlfName
Set-HyperlinkDisplay $wsZhCn 2 4 $newZhXlfName

$wsZhCn.Range("E2").Value = "2016-03-13 09:03:29"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = $newMdName
Set-HyperlinkDisplay $wsDeDe 2 1 $newMdName

$wsDeDe.Range("D2").Value = $newDeXlfName
Set-HyperlinkDisplay $wsDeDe 2 4 $newDeXlfName

$wsDeDe.Range("E2").Value = "2016-03-13 09:03:33"
